$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'96.780.99"
$ws.Range("E2").Value = "'  -1.16%  "

$ws.Range("D3").Value = "'3.337.88"
$ws.Range("E3").Value = "'  -2.61%  "

$ws.Range("E4").Value = "'  +0.04%  "

$ws.Range("D5").Value = "'250.21"

$ws.Range("D6").Value = "'655.71"
$ws.Range("E6").Value = "'  -0.16%  "

$ws.Range("E7").Value = "'  -5.33%  "

$ws.Range("D8").Value = "'0.422"
$ws.Range("E8").Value = "'  -2.51%  "

$ws.Range("D9").Value = "'1.00"
$ws.Range("E9").Value = "'  +0.06%  "

$ws.Range("E10").Value = "'  -5.73%  "

$ws.Range("D11").Value = "'3.334.79"
$ws.Range("E11").Value = "'  -2.65%  "

$ws.Range("D12").Value = "'0.206"
$ws.Range("E12").Value = "'  -2.74%  "

$ws.Range("D13").Value = "'40.68"
$ws.Range("E13").Value = "'  -3.18%  "

$ws.Range("D14").Value = "'96.547.11"
$ws.Range("E14").Value = "'  -1.02%  "

$ws.Range("D15").Value = "'6.09"
$ws.Range("E15").Value = "'  -3.73%  "

$ws.Range("E16").Value = "'  -2.86%  "

$ws.Range("D17").Value = "'3.958.56"
$ws.Range("E17").Value = "'  -2.75%  "

$ws.Range("D18").Value = "'8.71"
$ws.Range("E18").Value = "'  +0.73%  "

$ws.Range("D19").Value = "'3.374.85"
$ws.Range("E19").Value = "'  -1.50%  "

$ws.Range("D20").Value = "'0.570"
$ws.Range("E20").Value = "'  +13.38%  "

$ws.Range("D21").Value = "'17.43"
$ws.Range("E21").Value = "'  -0.92%  "

$ws.Range("D22").Value = "'10.67"
$ws.Range("E22").Value = "'  -0.96%  "

$ws.Range("D23").Value = "'508.35"
$ws.Range("E23").Value = "'  +0.46%  "

$ws.Range("D24").Value = "'3.33"
$ws.Range("E24").Value = "'  -4.08%  "

$ws.Range("E25").Value = "'  -4.05%  "

$ws.Range("E26").Value = "'  +6.24%  "

$ws.Range("D27").Value = "'96.52"
$ws.Range("E27").Value = "'  -2.65%  "

$ws.Range("D28").Value = "'12.10"
$ws.Range("E28").Value = "'  -5.61%  "

$ws.Range("B29").Value = "'Hedera"
$ws.Range("C29").Value = "'https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D29").Value = "'0.146"
$ws.Range("E29").Value = "'  -4.09%  "

$ws.Range("B30").Value = "'InternetComputer(DFINITY)"
$ws.Range("C30").Value = "'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D30").Value = "'11.41"
$ws.Range("E30").Value = "'  +0.05%  "

$ws.Range("B31").Value = "'Dai"
$ws.Range("C31").Value = "'https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D31").Value = "'1.00"
$ws.Range("E31").Value = "'  -0.16%  "

$ws.Range("B32").Value = "'Cronos"
$ws.Range("C32").Value = "'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D32").Value = "'0.188"
$ws.Range("E32").Value = "'  -7.05%  "

$ws.Range("B33").Value = "'PancakeSwap"
$ws.Range("C33").Value = "'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D33").Value = "'2.53"
$ws.Range("E33").Value = "'  +10.76%  "

$ws.Range("B34").Value = "'Binance-PegBSC-USD"
$ws.Range("C34").Value = "'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D34").Value = "'1.00"
$ws.Range("E34").Value = "'  +0.27%  "

$ws.Range("B35").Value = "'PolygonEcosystemToken"
$ws.Range("C35").Value = "'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D35").Value = "'0.553"
$ws.Range("E35").Value = "'  -3.99%  "

$ws.Range("B36").Value = "'EthereumClassic"
$ws.Range("C36").Value = "'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D36").Value = "'28.38"
$ws.Range("E36").Value = "'  -5.02%  "

$ws.Range("B37").Value = "'Fetch.AI"
$ws.Range("C37").Value = "'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D37").Value = "'1.51"
$ws.Range("E37").Value = "'  +4.52%  "

$ws.Range("B38").Value = "'RenderToken"
$ws.Range("C38").Value = "'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D38").Value = "'7.83"
$ws.Range("E38").Value = "'  -0.03%  "

$ws.Range("B39").Value = "'USDe"
$ws.Range("C39").Value = "'https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D39").Value = "'1.00"
$ws.Range("E39").Value = "'  +0.01%  "

$ws.Range("B40").Value = "'Kaspa"
$ws.Range("C40").Value = "'https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").Value = "'0.151"
$ws.Range("E40").Value = "'  -2.52%  "

$ws.Range("B41").Value = "'Bittensor"
$ws.Range("C41").Value = "'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D41").Value = "'506.58"
$ws.Range("E41").Value = "'  -3.17%  "

$ws.Range("B42").Value = "'WhiteBITCoin"
$ws.Range("C42").Value = "'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D42").Value = "'24.36"
$ws.Range("E42").Value = "'  -1.58%  "

$ws.Range("B43").Value = "'VeChain"
$ws.Range("C43").Value = "'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").Value = "'0.0434"
$ws.Range("E43").Value = "'  +3.54%  "

$ws.Range("B44").Value = "'MantraDAO"
$ws.Range("C44").Value = "'https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D44").Value = "'3.69"
$ws.Range("E44").Value = "'  -1.08%  "

$ws.Range("D45").Value = "'0.836"
$ws.Range("E45").Value = "'  -4.41%  "

$ws.Range("B46").Value = "'Filecoin"
$ws.Range("C46").Value = "'https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D46").Value = "'5.57"
$ws.Range("E46").Value = "'  -0.23%  "

$ws.Range("D47").Value = "'1.67"
$ws.Range("E47").Value = "'  +4.31%  "

$ws.Range("B48").Value = "'Cosmos"
$ws.Range("C48").Value = "'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D48").Value = "'8.49"
$ws.Range("E48").Value = "'  +2.73%  "

$ws.Range("B49").Value = "'OKB"
$ws.Range("C49").Value = "'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D49").Value = "'54.65"
$ws.Range("E49").Value = "'  +6.60%  "

$ws.Range("B50").Value = "'dogwifhat"
$ws.Range("C50").Value = "'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D50").Value = "'3.10"
$ws.Range("E50").Value = "'  -6.96%  "

$ws.Range("B51").Value = "'Monero"
$ws.Range("C51").Value = "'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D51").Value = "'162.06"
$ws.Range("E51").Value = "'  +0.33%  "

